# Weekly fruit/hortaliza refresh: a new "Segunda" quality record for
# Terminal Hortofrutícola Agro Chillán - Brócoli was collected for the
# week of 2022-09-05. It is inserted as the new first data row (row 229),
# pushing the existing historical rows (229-307) down by one (230-308).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 229, shifting rows 229:307 down to 230:308.
$ws.Rows("229:229").Insert()

# Populate the newly inserted row with the latest week's reading.
$ws.Range("A229").Value = 7
$ws.Range("B229").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C229").Value = "Ñuble"
$ws.Range("D229").Value = Get-Date -Year 2022 -Month 9 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("E229").Value = 16
$ws.Range("F229").Value = 100112023
$ws.Range("G229").Value = "Brócoli"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Segunda"
$ws.Range("J229").Value = 150
$ws.Range("K229").Value = 800
$ws.Range("L229").Value = 800
$ws.Range("M229").Value = 800
$ws.Range("N229").Value = "`$/unidad"
$ws.Range("O229").Value = "Región del Maule"
$ws.Range("P229").Value = 800
$ws.Range("Q229").Value = 1
$ws.Range("R229").Value = "Hortaliza"
